# Rename the "discretization_*" headers to the shorter "disc_*" form, and
# refresh the discretized/validation values below them (the decision-tree
# classifier now buckets "grade" into 3 categories, plus the assorted
# value corrections noted in the commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:P1) -------------------------------------------------
$headers = @(
    "disc_sqft_lot15",
    "disc_sqft_living15",
    "disc_yr_renovated",
    "disc_yr_built",
    "disc_sqft_basement",
    "disc_sqft_above",
    "disc_grade",
    "disc_condition",
    "disc_view",
    "disc_floors",
    "disc_bedrooms",
    "disc_bathrooms",
    "disc_sqft_lot",
    "disc_sqft_living",
    "disc_lat_long",
    "caro"
)
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# --- Data rows (A2:P11) --------------------------------------------------
$data = @(
    @(1,1,0,1,1,0,2,1,0,1,1,0,1,1,0,1),
    @(1,1,0,1,1,0,2,1,0,1,1,0,1,1,0,1),
    @(1,1,0,1,1,0,2,1,0,1,1,0,1,1,0,0),
    @(1,1,0,1,1,0,2,1,0,1,0,0,1,1,0,1),
    @(1,1,0,1,1,0,2,1,0,1,1,0,1,1,0,0),
    @(1,0,0,1,1,1,2,1,0,1,0,0,1,0,0,1),
    @(1,0,0,1,1,1,2,1,0,1,0,0,1,0,0,2),
    @(1,1,0,1,1,0,1,1,0,1,1,0,1,1,0,0),
    @(1,1,0,1,1,0,1,1,0,1,1,0,1,1,0,1),
    @(1,1,0,1,1,0,1,1,0,1,1,0,1,1,0,0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $i + 2
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($rowIndex, $c + 1).Value = $rowData[$c]
    }
}
